$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "modified" timestamp
$ws.Range("B21").Value = "2023-08-17T07:59:39+00:00"

# Shift terms: row24 B becomes "subject", row25 B becomes "variable",
# row26 B becomes "community maturity level", and rows 27-29 become the
# three maturity-level terms with C pointing to "community maturity level".
$ws.Range("B24").Value = "subject"
$ws.Range("C24").Value = ""

$ws.Range("B25").Value = "variable"
$ws.Range("C25").Value = ""

$ws.Range("B26").Value = "community maturity level"
$ws.Range("C26").Value = ""

$ws.Range("B27").Value = "emerging"
$ws.Range("C27").Value = "community maturity level"

$ws.Range("B28").Value = "developing"
$ws.Range("C28").Value = "community maturity level"

$ws.Range("B29").Value = "mature"
$ws.Range("C29").Value = "community maturity level"

# Delete rows 30-44 (old terms no longer present in the new vocabulary)
$rng = $ws.Range("A30:AM44")
$rng.Delete()
